# Add "Save" column (H) to the s_vals sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, styled like the other header cells (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells H2:H3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
